# Add survey name on importing new survey responses:
# Rename the first sheet from "Basic Clinic Data..." to "Test Survey",
# and update which sheet/cell is active & selected in each sheet view.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Give the survey-response sheet a proper name.
$ws1.Name = "Test Survey"

# Make sure the second sheet's selection is updated first...
$ws2.Activate()
$ws2.Range("G7").Select()

# ...then make the renamed survey sheet the active/selected tab,
# with its own new selection.
$ws1.Activate()
$ws1.Range("E1").Select()
